$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item("user details")

# Rename the second sheet
$ws2.Name = "user_details"

# Populate the user_details sheet with header row + data row
$ws2.Range("A1").Value = "lastname"
$ws2.Range("B1").Value = "usergroup"
$ws2.Range("C1").Value = "userrole"

$ws2.Range("A2").Value = " ln"
$ws2.Range("B2").Value = "Convirza-Live"
$ws2.Range("C2").Value = "Admin"

# Update selections/active sheet to match the new state
$ws1.Range("A2").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("M15").Select() | Out-Null
